$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: insert two new "Unnamed" columns before the existing
# "Dia da semana" header, shifting that label from D1 to F1 ---
$origHeader = $ws.Range("D1").Text
$ws.Range("F1").Value = $origHeader
$ws.Range("D1").Value = "Unnamed: 3"
$ws.Range("E1").Value = "Unnamed: 4"

# match the bold / centered / bordered look of the other header cells
$ws.Range("D1:F1").Font.Bold = $true
$ws.Range("D1:F1").HorizontalAlignment = -4108
$ws.Range("D1:F1").VerticalAlignment = -4160
$ws.Range("D1:F1").Borders.LineStyle = 1

# --- Data rows: the weekday name that used to live in column D now
# belongs in column F; columns D and E become blank placeholder columns ---
for ($i = 2; $i -le 369; $i++) {
    $dCell = "D" + $i
    $eCell = "E" + $i
    $fCell = "F" + $i
    $weekday = $ws.Range($dCell).Text
    $ws.Range($fCell).Value = $weekday
    $ws.Range($dCell).Value = ""
    $ws.Range($eCell).Value = ""
}

# --- Two new rows of data appended at the bottom ---
$ws.Range("A370").NumberFormat = "@"
$ws.Range("A370").Value = "03/16/2021"
$ws.Range("B370").Value = 0.44
$ws.Range("C370").Value = 0.43
$ws.Range("F370").Value = "Terça-feira"

$ws.Range("A371").NumberFormat = "@"
$ws.Range("A371").Value = "03/17/2021"
$ws.Range("B371").Value = 0.43
$ws.Range("C371").Value = 0.42
$ws.Range("F371").Value = "Quarta-feira"
